# Replace the legacy field-code construct (fldChar begin / instrText /
# fldChar end) that implements the M2Doc "sampleTable" call with a plain
# literal-text run containing the equivalent curly-brace token, as used
# by the new TokenIteratorFieldRewriterSplit parser.
$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$cell = $table.Rows.Item(2).Cells.Item(2)
$paragraph = $cell.Range.Paragraphs.Item(1)

$newParagraphXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>{m:'A sample table'.sampleTable()}</w:t></w:r></w:p>"

$null = $paragraph.Range.InsertXML($newParagraphXml)
